$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column AR (27-jul) ---
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Range("AR1").Value = "27-jul"
# match the formatting used by the rest of the header row (bold, thin border,
# centered horizontally, top-aligned vertically)
$wsPrix.Range("AR1").Font.Bold = $true
$wsPrix.Range("AR1").HorizontalAlignment = -4108
$wsPrix.Range("AR1").VerticalAlignment = -4160
$wsPrix.Range("AR1").Borders.LineStyle = 1
$wsPrix.Range("AR1").Borders.Weight = 2

$arValues = @{
    2  = 60.5
    3  = 49.86
    4  = 48.9
    5  = 32.09
    6  = 37.92
    7  = 40.01
    8  = 40.16
    9  = 50
    10 = 32.34
    11 = 23.69
    12 = 14.77
    13 = 33.17
    14 = 31.66
    15 = 15.04
    16 = 8.85
    17 = 12.93
    18 = 16.14
    19 = 24
    20 = 23.3
    21 = 35.32
    22 = 40
    23 = 62.67
    24 = 88.64
    25 = 66.95
}

foreach ($row in $arValues.Keys) {
    $wsPrix.Cells.Item($row, 44).Value = $arValues[$row]
}

# --- Sheet "Gaz": add row 41 (2025-07-25, 31.7) ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A41").NumberFormat = "@"
$wsGaz.Range("A41").Value = "2025-07-25"
$wsGaz.Range("B41").Value = 31.7

# --- Sheet "CO2": add row 41 (2025-07-25, 70.7) ---
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A41").NumberFormat = "@"
$wsCo2.Range("A41").Value = "2025-07-25"
$wsCo2.Range("B41").Value = 70.7
